# Auto-generated Excel COM-interop script
# Applies scheduled market-price data refresh to the Leve profit sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1744
$ws.Range("J43").Value = 1743.75
$ws.Range("L43").Value = 1743.75
$ws.Range("N43").Value = -1881.75
$ws.Range("H121").Value = 1639.9375
$ws.Range("J121").Value = 1774.2142
$ws.Range("L121").Value = 5322.642599999999
$ws.Range("N121").Value = -8816.642599999999
$ws.Range("H132").Value = 3029.889
$ws.Range("I132").Value = 2938.72
$ws.Range("J132").Value = 4169.5
$ws.Range("K132").Value = 8816.16
$ws.Range("L132").Value = 12508.5
$ws.Range("M132").Value = -6286.16
$ws.Range("N132").Value = -17568.5
$ws.Range("H141").Value = 4600
$ws.Range("I141").Value = 4400
$ws.Range("J141").Value = 4666.6665
$ws.Range("K141").Value = 13200
$ws.Range("L141").Value = 13999.9995
$ws.Range("M141").Value = -8020
$ws.Range("N141").Value = -24359.9995

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3070.4324
$ws.Range("I32").Value = 2560.1714
$ws.Range("J32").Value = 12000
$ws.Range("K32").Value = 2560.1714
$ws.Range("L32").Value = 12000
$ws.Range("M32").Value = -2273.1714
$ws.Range("N32").Value = -12574
$ws.Range("H45").Value = 3104.8372
$ws.Range("I45").Value = 2624.7
$ws.Range("K45").Value = 2624.7
$ws.Range("M45").Value = -2247.7

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H27").Value = 35180.668
$ws.Range("J27").Value = 35180.668
$ws.Range("L27").Value = 35180.668
$ws.Range("N27").Value = -35564.668
$ws.Range("H100").Value = 27610.75
$ws.Range("J100").Value = 27610.75
$ws.Range("L100").Value = 27610.75
$ws.Range("N100").Value = -29774.75
$ws.Range("H105").Value = 3104.842
$ws.Range("I105").Value = 3659.111
$ws.Range("J105").Value = 2606
$ws.Range("K105").Value = 3659.111
$ws.Range("L105").Value = 2606
$ws.Range("M105").Value = -1912.111
$ws.Range("N105").Value = -6100

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3282.2727
$ws.Range("I31").Value = 3005.3076
$ws.Range("J31").Value = 3462.3
$ws.Range("K31").Value = 3005.3076
$ws.Range("L31").Value = 3462.3
$ws.Range("M31").Value = -2710.3076
$ws.Range("N31").Value = -4052.3
$ws.Range("H34").Value = 3282.2727
$ws.Range("I34").Value = 3005.3076
$ws.Range("J34").Value = 3462.3
$ws.Range("K34").Value = 3005.3076
$ws.Range("L34").Value = 3462.3
$ws.Range("M34").Value = -2803.3076
$ws.Range("N34").Value = -3866.3
$ws.Range("H58").Value = 23050.217
$ws.Range("I58").Value = 1514.2307
$ws.Range("K58").Value = 1514.2307
$ws.Range("M58").Value = -1311.2307
$ws.Range("H96").Value = 3748.6
$ws.Range("J96").Value = 3748.6
$ws.Range("L96").Value = 3748.6
$ws.Range("N96").Value = -9240.6
$ws.Range("H122").Value = 2627.2856
$ws.Range("I122").Value = 2627.2856
$ws.Range("K122").Value = 7881.8568
$ws.Range("M122").Value = -5431.8568
$ws.Range("H132").Value = 3098.72
$ws.Range("I132").Value = 2087.1
$ws.Range("J132").Value = 7145.2
$ws.Range("K132").Value = 6261.299999999999
$ws.Range("L132").Value = 21435.6
$ws.Range("M132").Value = -3731.299999999999
$ws.Range("N132").Value = -26495.6
$ws.Range("H136").Value = 23050.217
$ws.Range("I136").Value = 1514.2307
$ws.Range("K136").Value = 4542.6921
$ws.Range("M136").Value = -1992.6921

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 2471.2683
$ws.Range("I2").Value = 3057.606
$ws.Range("J2").Value = 52.625
$ws.Range("K2").Value = 18345.636
$ws.Range("L2").Value = 315.75
$ws.Range("M2").Value = -18232.636
$ws.Range("N2").Value = -541.75
$ws.Range("H5").Value = 1365.5385
$ws.Range("I5").Value = 1229.3334
$ws.Range("K5").Value = 3688.0002
$ws.Range("M5").Value = -3576.0002
$ws.Range("H17").Value = 404.8
$ws.Range("I17").Value = 133.16667
$ws.Range("J17").Value = 812.25
$ws.Range("K17").Value = 399.50001
$ws.Range("L17").Value = 2436.75
$ws.Range("M17").Value = -230.50001
$ws.Range("N17").Value = -2774.75
$ws.Range("H34").Value = 836.1539
$ws.Range("I34").Value = 510
$ws.Range("J34").Value = 895.4545000000001
$ws.Range("K34").Value = 1530
$ws.Range("L34").Value = 2686.3635
$ws.Range("M34").Value = -1446
$ws.Range("N34").Value = -2854.3635
$ws.Range("H39").Value = 1316.1333
$ws.Range("J39").Value = 1274.4286
$ws.Range("L39").Value = 3823.2858
$ws.Range("N39").Value = -4411.2858
$ws.Range("H55").Value = 2615.125
$ws.Range("J55").Value = 2615.125
$ws.Range("L55").Value = 7845.375
$ws.Range("N55").Value = -8199.375
$ws.Range("H122").Value = 651.4286
$ws.Range("I122").Value = 427
$ws.Range("J122").Value = 1998
$ws.Range("K122").Value = 3843
$ws.Range("L122").Value = 17982
$ws.Range("M122").Value = -1393
$ws.Range("N122").Value = -22882
$ws.Range("H129").Value = 1453.4445
$ws.Range("I129").Value = 510.125
$ws.Range("K129").Value = 1530.375
$ws.Range("M129").Value = 3469.625
$ws.Range("H131").Value = 781.9794000000001
$ws.Range("I131").Value = 430
$ws.Range("J131").Value = 797.1183
$ws.Range("K131").Value = 1290
$ws.Range("L131").Value = 2391.3549
$ws.Range("M131").Value = 3750
$ws.Range("N131").Value = -12471.3549
$ws.Range("H132").Value = 1018.625
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("H135").Value = 1365.5385
$ws.Range("I135").Value = 1229.3334
$ws.Range("K135").Value = 11064.0006
$ws.Range("M135").Value = -8529.000599999999

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 46500
$ws.Range("J108").Value = 46500
$ws.Range("L108").Value = 46500
$ws.Range("N108").Value = -54180
$ws.Range("H113").Value = 3680
$ws.Range("I113").Value = 3000
$ws.Range("J113").Value = 4020
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 4020
$ws.Range("M113").Value = -830
$ws.Range("N113").Value = -8360
$ws.Range("H126").Value = 4503.8965
$ws.Range("I126").Value = 3215.7368
$ws.Range("J126").Value = 6951.4
$ws.Range("K126").Value = 9647.2104
$ws.Range("L126").Value = 20854.2
$ws.Range("M126").Value = -7177.2104
$ws.Range("N126").Value = -25794.2

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5725.5
$ws.Range("I22").Value = 5725.5
$ws.Range("K22").Value = 5725.5
$ws.Range("M22").Value = -5430.5
$ws.Range("H27").Value = 5725.5
$ws.Range("I27").Value = 5725.5
$ws.Range("K27").Value = 5725.5
$ws.Range("M27").Value = -5618.5
$ws.Range("H40").Value = 6975.8335
$ws.Range("I40").Value = 3212.5
$ws.Range("J40").Value = 14502.5
$ws.Range("K40").Value = 3212.5
$ws.Range("L40").Value = 14502.5
$ws.Range("M40").Value = -3076.5
$ws.Range("N40").Value = -14774.5
$ws.Range("H68").Value = 4464.467
$ws.Range("I68").Value = 2163.3333
$ws.Range("K68").Value = 2163.3333
$ws.Range("M68").Value = -1414.3333
$ws.Range("H71").Value = 4464.467
$ws.Range("I71").Value = 2163.3333
$ws.Range("K71").Value = 10816.6665
$ws.Range("M71").Value = -7072.666499999999
$ws.Range("H82").Value = 2915
$ws.Range("I82").Value = 3034
$ws.Range("J82").Value = 2724.6
$ws.Range("K82").Value = 3034
$ws.Range("L82").Value = 2724.6
$ws.Range("M82").Value = -2673
$ws.Range("N82").Value = -3446.6
$ws.Range("H85").Value = 2915
$ws.Range("I85").Value = 3034
$ws.Range("J85").Value = 2724.6
$ws.Range("K85").Value = 3034
$ws.Range("L85").Value = 2724.6
$ws.Range("M85").Value = -1786
$ws.Range("N85").Value = -5220.6

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2130.875
$ws.Range("I122").Value = 1939.2858
$ws.Range("K122").Value = 5817.857400000001
$ws.Range("M122").Value = -3367.857400000001
$ws.Range("H136").Value = 17858340
$ws.Range("I136").Value = 26316866
$ws.Range("J136").Value = 1449.6111
$ws.Range("K136").Value = 78950598
$ws.Range("L136").Value = 4348.8333
$ws.Range("M136").Value = -78948048
$ws.Range("N136").Value = -9448.8333
